# Post Bid create api
# Adds a new "Bid" service row (row 37) to the ServicesList sheet, marks the
# prior row (Image Upload, row 36) "Coding" column as Done, and adjusts the
# sheet view to freeze the header row / scroll the visible area.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 36: mark "Coding" (column K) as Done -------------------------------
$ws.Range("K36").Value = "Done"

# --- Row 37: new "Bid" service entry ----------------------------------------
# Copy row 36's cell formatting down into row 37 first so borders/fills/
# number formats match the rest of the table.
$ws.Range("B36:O36").Copy()
$ws.Range("B37").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# Text / shared-string values (order mirrors how the new shared strings end
# up appended to the workbook's string table).
$ws.Range("B37").Value = "Bid"
$ws.Range("D37").Value = "WS-BID-01"
$ws.Range("C37").Value = "Create Post bid"
$ws.Range("E37").Value = "app.bid.post.insert"
$ws.Range("G37").Value = "bid"
$ws.Range("H37").Value = "/create"
$ws.Range("I37").Value = "POST"
$ws.Range("K37").Value = "Done"

# logActivity (column F) must be the literal text "true" (matching the rest
# of the column, which stores "true"/"false" as text, not booleans). Typing
# "true" straight into the cell gets auto-coerced to a Boolean, so round it
# through a text formula on a scratch cell and paste just the value back.
$ws.Range("Z1").Formula = '="true"'
$ws.Range("Z1").Copy()
$ws.Range("F37").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("Z1").ClearContents()

# Generated-code helper formulas, matching the pattern used by every other
# row in the table.
$ws.Range("M37").Formula = '=_xlfn.CONCAT("INSERT INTO ",CHAR(34),"M_CTL_CONFIG",CHAR(34)," VALUES(''",D37,"'',''CONNON_CONFIG'', 0, ''",C37,"'', ''{}'', 0, 0, CURRENT_TIMESTAMP, ''ATUL'', null, null);")'
$ws.Range("N37").Formula = '=_xlfn.CONCAT(IF(I37="GET","@GetMapping(",IF(I37="POST","@PostMapping(",IF(I37="DELETE","@DeleteMapping(",IF(I37="PUT","@PutMapping(","")))),CHAR(34),H37,CHAR(34),")")'
$ws.Range("O37").Formula = '=_xlfn.CONCAT("@ServiceInfo(serviceCode = ",CHAR(34),D37,,CHAR(34),", serviceName = ",CHAR(34),C37,CHAR(34), ", queryId = ",CHAR(34),E37,CHAR(34),", logActivity =",F37,")")'

# --- Sheet view: freeze the header row and scroll / select -----------------
$ws.Range("A4").Select()
$excel.ActiveWindow.FreezePanes = $true
$excel.ActiveWindow.ScrollRow = 18
$ws.Range("H41").Select()
